$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("A19").Value = 112392912
$ws.Range("B19").Value = 90800
$ws.Range("C19").Value = 'Ovaliderad'
$ws.Range("D19").Value = 'LC'
$ws.Range("E19").Value = 4364
$ws.Range("F19").Value = 'Dropptaggsvamp'
$ws.Range("G19").Value = 'Hydnellum ferrugineum'
$ws.Range("H19").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("P19").Value = 'Järbäcken, Ög'
$ws.Range("Q19").Value = 562788
$ws.Range("R19").Value = 6504516
$ws.Range("S19").Value = 50
$ws.Range("T19").Value = 'Östergötland'
$ws.Range("U19").Value = 'Norrköping'
$ws.Range("V19").Value = 'Östergötland'
$ws.Range("W19").Value = 'Kvillinge'
$ws.Range("Y19").Value = "'2023-09-29"
$ws.Range("Z19").Value = '13:58'
$ws.Range("AA19").Value = "'2023-09-29"
$ws.Range("AB19").Value = '13:58'
$ws.Range("AD19").Value = $false
$ws.Range("AE19").Value = $false
$ws.Range("AG19").Value = $false
$ws.Range("AW19").Value = 'Jens Johannesson'
$ws.Range("AX19").Value = 'Jens Johannesson'

# Row 20
$ws.Range("A20").Value = 112392668
$ws.Range("B20").Value = 90155
$ws.Range("C20").Value = 'Ovaliderad'
$ws.Range("D20").Value = 'LC'
$ws.Range("E20").Value = 6031
$ws.Range("F20").Value = 'Blomkålssvamp'
$ws.Range("G20").Value = 'Sparassis crispa'
$ws.Range("H20").Value = '(Wulfen:Fr.) Fr.'
$ws.Range("P20").Value = 'Järbäcken, Ög'
$ws.Range("Q20").Value = 562710
$ws.Range("R20").Value = 6504599
$ws.Range("S20").Value = 50
$ws.Range("T20").Value = 'Östergötland'
$ws.Range("U20").Value = 'Norrköping'
$ws.Range("V20").Value = 'Östergötland'
$ws.Range("W20").Value = 'Kvillinge'
$ws.Range("Y20").Value = "'2023-09-29"
$ws.Range("Z20").Value = '13:32'
$ws.Range("AA20").Value = "'2023-09-29"
$ws.Range("AB20").Value = '13:32'
$ws.Range("AD20").Value = $false
$ws.Range("AE20").Value = $false
$ws.Range("AG20").Value = $false
$ws.Range("AW20").Value = 'Jens Johannesson'
$ws.Range("AX20").Value = 'Jens Johannesson'

# Row 21
$ws.Range("A21").Value = 112392511
$ws.Range("B21").Value = 89936
$ws.Range("C21").Value = 'Ovaliderad'
$ws.Range("D21").Value = 'LC'
$ws.Range("E21").Value = 5420
$ws.Range("F21").Value = 'Grovticka'
$ws.Range("G21").Value = 'Phaeolus schweinitzii'
$ws.Range("H21").Value = '(Fr.) Pat.'
$ws.Range("P21").Value = 'Järbäcken, Ög'
$ws.Range("Q21").Value = 562765
$ws.Range("R21").Value = 6504441
$ws.Range("S21").Value = 50
$ws.Range("T21").Value = 'Östergötland'
$ws.Range("U21").Value = 'Norrköping'
$ws.Range("V21").Value = 'Östergötland'
$ws.Range("W21").Value = 'Kvillinge'
$ws.Range("Y21").Value = "'2023-09-29"
$ws.Range("Z21").Value = '13:32'
$ws.Range("AA21").Value = "'2023-09-29"
$ws.Range("AB21").Value = '13:32'
$ws.Range("AC21").Value = 'Gammal tallskog.'
$ws.Range("AD21").Value = $false
$ws.Range("AE21").Value = $false
$ws.Range("AG21").Value = $false
$ws.Range("AW21").Value = 'Jens Johannesson'
$ws.Range("AX21").Value = 'Jens Johannesson'

# Row 22
$ws.Range("A22").Value = 112392758
$ws.Range("B22").Value = 93539
$ws.Range("C22").Value = 'Ovaliderad'
$ws.Range("D22").Value = 'LC'
$ws.Range("E22").Value = 2180
$ws.Range("F22").Value = 'Blåmossa'
$ws.Range("G22").Value = 'Leucobryum glaucum'
$ws.Range("H22").Value = '(Hedw.) Ångstr.'
$ws.Range("P22").Value = 'Järbäcken, Ög'
$ws.Range("Q22").Value = 562788
$ws.Range("R22").Value = 6504516
$ws.Range("S22").Value = 50
$ws.Range("T22").Value = 'Östergötland'
$ws.Range("U22").Value = 'Norrköping'
$ws.Range("V22").Value = 'Östergötland'
$ws.Range("W22").Value = 'Kvillinge'
$ws.Range("Y22").Value = "'2023-09-29"
$ws.Range("Z22").Value = '13:58'
$ws.Range("AA22").Value = "'2023-09-29"
$ws.Range("AB22").Value = '13:58'
$ws.Range("AC22").Value = 'Blåmossa i större bestånd.'
$ws.Range("AD22").Value = $false
$ws.Range("AE22").Value = $false
$ws.Range("AG22").Value = $false
$ws.Range("AW22").Value = 'Jens Johannesson'
$ws.Range("AX22").Value = 'Jens Johannesson'

Write-Host "Added rows 19-22"